$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Unique" id values in column A for rows 2-4
$ws.Range("A2").Value = "9042"
$ws.Range("A3").Value = "9035"
$ws.Range("A4").Value = "9001"

# Delete row 5 (previously Unique=1002, Date=45064, Comments="test gsdgdfbfdbvdf dsv")
$ws.Rows.Item(5).Delete()

# Update the selection to match the target workbook state
$ws.Range("A5").Select()
